$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (rows 2 through 319) from serial date 45177 to 45178
for ($row = 2; $row -le 319; $row++) {
    $ws.Cells.Item($row, 3).Value = 45178
}
